# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) timestamps on row 4 of the
# zh-cn and de-de report sheets, to reflect the regenerated report run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-25 13:34:24"
$wsZhCn.Range("G4").Value = "2016-01-25 13:35:10"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-25 13:34:33"
$wsDeDe.Range("G4").Value = "2016-01-25 13:35:28"
